# LOOKUPTABLE change for resistor rules
# Adds a new "RESISTOR_PREFIX" lookup sheet as the first sheet of the
# workbook, ahead of the existing TABLE and ROUTING_RULES sheets.

$wb = $excel.ActiveWorkbook

# Insert a brand new worksheet before the current first sheet (TABLE) so
# that it becomes sheet #1, pushing TABLE/ROUTING_RULES down.
$ws = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$ws.Name = "RESISTOR_PREFIX"

# Header row
$ws.Range("A1").Value = "Prefix"
$ws.Range("B1").Value = "Rating_Value"
$ws.Range("C1").Value = "Rating_Unit"
$ws.Range("D1").Value = "Vendor"
$ws.Range("E1").Value = "Priority"

# Data rows
$ws.Range("A2").Value = "WR02X"
$ws.Range("B2").Value = 0.1
$ws.Range("C2").Value = "W"
$ws.Range("D2").Value = "WALSIN"
$ws.Range("E2").Value = 1

$ws.Range("A3").Value = "WR04X"
$ws.Range("B3").Value = 0.125
$ws.Range("C3").Value = "W"
$ws.Range("D3").Value = "WALSIN"
$ws.Range("E3").Value = 1

$ws.Range("A4").Value = "WR06X"
$ws.Range("B4").Value = 0.25
$ws.Range("C4").Value = "W"
$ws.Range("D4").Value = "WALSIN"
$ws.Range("E4").Value = 1

$ws.Range("A5").Value = "CRCW0"
$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = "W"
$ws.Range("D5").Value = "VISHAY"
$ws.Range("E5").Value = 2

$ws.Range("A6").Value = "ERJ2G"
$ws.Range("B6").Value = 0.1
$ws.Range("C6").Value = "W"
$ws.Range("D6").Value = "PANASONIC"
$ws.Range("E6").Value = 2

# Make the new sheet the active/selected tab, and restore the scroll
# position / selection on the TABLE sheet to match its prior state.
$tableWs = $wb.Worksheets.Item("TABLE")
$tableWs.Activate()
$tableWs.Range("G154").Select()

$ws.Activate()
$ws.Range("A1").Select()
